$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '95.283.03'
Set-TextValue 'E2' '  -2.26%  '
Set-TextValue 'D3' '3.596.67'
Set-TextValue 'E3' '  -3.40%  '
Set-TextValue 'D4' '2.66'
Set-TextValue 'E4' '  +38.15%  '
Set-TextValue 'D5' '0.999'
Set-TextValue 'E5' '  -0.16%  '
Set-TextValue 'D6' '222.58'
Set-TextValue 'E6' '  -6.20%  '
Set-TextValue 'D7' '636.19'
Set-TextValue 'E7' '  -3.53%  '
Set-TextValue 'D8' '0.417'
Set-TextValue 'E8' '  -4.72%  '
Set-TextValue 'D9' '1.18'
Set-TextValue 'E9' '  +10.43%  '
Set-TextValue 'D10' '0.999'
Set-TextValue 'E10' '  -0.03%  '
Set-TextValue 'D11' '3.593.21'
Set-TextValue 'E11' '  -3.42%  '
Set-TextValue 'D12' '47.66'
Set-TextValue 'E12' '  +6.16%  '
Set-TextValue 'E13' '  +1.00%  '
Set-TextValue 'D14' '0.0000290'
Set-TextValue 'E14' '  -9.01%  '
Set-TextValue 'D15' '6.47'
Set-TextValue 'E15' '  -6.71%  '
Set-TextValue 'D16' '4.293.47'
Set-TextValue 'E16' '  -2.81%  '
Set-TextValue 'D17' '95.013.08'
Set-TextValue 'E17' '  -2.11%  '
Set-TextValue 'D18' '22.69'
Set-TextValue 'E18' '  +20.53%  '
Set-TextValue 'D19' '8.87'
Set-TextValue 'E19' '  -1.80%  '
Set-TextValue 'D20' '13.89'
Set-TextValue 'E20' '  +6.26%  '
Set-TextValue 'D21' '3.593.65'
Set-TextValue 'E21' '  -3.29%  '
Set-TextValue 'D22' '0.537'
Set-TextValue 'E22' '  +5.49%  '
Set-TextValue 'D23' '0.280'
Set-TextValue 'E23' '  +45.99%  '
Set-TextValue 'D24' '512.03'
Set-TextValue 'E24' '  -2.74%  '
Set-TextValue 'D25' '3.23'
Set-TextValue 'E25' '  -7.35%  '
Set-TextValue 'D26' '120.40'
Set-TextValue 'E26' '  +12.89%  '
Set-TextValue 'D27' '0.0000201'
Set-TextValue 'E27' '  -11.13%  '
Set-TextValue 'D28' '6.78'
Set-TextValue 'E28' '  -1.86%  '
Set-TextValue 'D29' '3.777.55'
Set-TextValue 'E29' '  -3.63%  '
Set-TextValue 'D30' '12.65'
Set-TextValue 'E30' '  -6.88%  '
Set-TextValue 'D31' '12.73'
Set-TextValue 'E31' '  +0.21%  '
Set-TextValue 'D32' '3.01'
Set-TextValue 'E32' '  -0.50%  '
Set-TextValue 'E33' '  +0.18%  '
Set-TextValue 'E34' '  +0.12%  '
Set-TextValue 'B35' 'PolygonEcosystemToken'
Set-TextValue 'C35' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D35' '0.614'
Set-TextValue 'E35' '  +3.34%  '
Set-TextValue 'D36' '32.45'
Set-TextValue 'E36' '  -0.72%  '
Set-TextValue 'B37' 'Cronos'
Set-TextValue 'C37' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D37' '0.179'
Set-TextValue 'E37' '  -7.07%  '
Set-TextValue 'D38' '1.74'
Set-TextValue 'E38' '  -5.32%  '
Set-TextValue 'D40' '8.29'
Set-TextValue 'E40' '  -5.60%  '
Set-TextValue 'B41' 'Filecoin'
Set-TextValue 'C41' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D41' '7.06'
Set-TextValue 'E41' '  +4.60%  '
Set-TextValue 'B42' 'Bittensor'
Set-TextValue 'C42' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D42' '575.35'
Set-TextValue 'E42' '  -10.54%  '
Set-TextValue 'D43' '41.96'
Set-TextValue 'E43' '  +3.04%  '
Set-TextValue 'D44' '0.503'
Set-TextValue 'E44' '  +1.62%  '
Set-TextValue 'D45' '0.0507'
Set-TextValue 'E45' '  +11.04%  '
Set-TextValue 'D46' '0.155'
Set-TextValue 'E46' '  -6.60%  '
Set-TextValue 'D47' '0.955'
Set-TextValue 'E47' '  -1.74%  '
Set-TextValue 'D48' '1.93'
Set-TextValue 'E48' '  -4.83%  '
Set-TextValue 'D49' '8.90'
Set-TextValue 'E49' '  +2.57%  '
Set-TextValue 'D50' '229.07'
Set-TextValue 'E50' '  +10.24%  '
Set-TextValue 'D51' '23.48'
Set-TextValue 'E51' '  -0.70%  '
